# Generate Report for Handoff
# Updates the "07bb04c5-..." handoff-related rows (7,8,9,10,12,13) across the
# Overview, zh-cn and de-de sheets:
#   - Overview!G  (Latest HO Xliff Generate Date)  -> 2016-08-26 14:31:27
#   - de-de!H     (Latest Handoff Datetime, de-de)  -> 2016-08-26 14:31:27 (was in sync with Overview!G)
#   - zh-cn!H     (Latest Handoff Datetime, zh-cn)  -> 2016-08-26 14:31:21
#   - zh-cn!E     (Priority)                        -> ht
#   - de-de!E     (Priority)                        -> ht

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 13)

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-26 14:31:27"
    $wsDeDe.Range("H$r").Value     = "2016-08-26 14:31:27"
    $wsZhCn.Range("H$r").Value     = "2016-08-26 14:31:21"
    $wsZhCn.Range("E$r").Value     = "ht"
    $wsDeDe.Range("E$r").Value     = "ht"
}
